# Insert a new data row at row 237 (pushing existing rows 237-278 down to 238-279)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 237..278 down by one, creating a new blank row 237.
$ws.Rows.Item(237).Insert()

# Fill in the new row with the new weekly record.
$ws.Cells.Item(237, 1).Value  = 11
$ws.Cells.Item(237, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(237, 3).Value  = "Bíobío"
$ws.Cells.Item(237, 4).Value  = 45209
$ws.Cells.Item(237, 5).Value  = 8
$ws.Cells.Item(237, 6).Value  = 100112032
$ws.Cells.Item(237, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(237, 8).Value  = "Sin especificar"
$ws.Cells.Item(237, 9).Value  = "Primera"
$ws.Cells.Item(237, 10).Value = 120
$ws.Cells.Item(237, 11).Value = 15000
$ws.Cells.Item(237, 12).Value = 15000
$ws.Cells.Item(237, 13).Value = 15000
$ws.Cells.Item(237, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(237, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(237, 16).Value = 300
$ws.Cells.Item(237, 17).Value = 50
$ws.Cells.Item(237, 18).Value = "Hortaliza"
